$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the values that currently live in rows 19 and 22 (before any rows are removed)
$b19 = $ws.Range("B19").Value()
$c19 = $ws.Range("C19").Value()
$d19 = $ws.Range("D19").Value()
$e19 = $ws.Range("E19").Value()
$f19 = $ws.Range("F19").Value()

$b22 = $ws.Range("B22").Value()
$c22 = $ws.Range("C22").Value()
$d22 = $ws.Range("D22").Value()
$e22 = $ws.Range("E22").Value()
$f22 = $ws.Range("F22").Value()

# Delete rows 13 through 22 (the rows that disappear from the final report)
$ws.Range("A13:F22").EntireRow.Delete()

# Row 11 becomes what used to be row 19
$ws.Range("B11").Value = $b19
$ws.Range("C11").Value = $c19
$ws.Range("D11").Value = $d19
$ws.Range("E11").Value = $e19
$ws.Range("F11").Value = $f19

# Row 12 becomes what used to be row 22
$ws.Range("B12").Value = $b22
$ws.Range("C12").Value = $c22
$ws.Range("D12").Value = $d22
$ws.Range("E12").Value = $e22
$ws.Range("F12").Value = $f22

# Restore the selection Excel left behind after the edit
$ws.Range("G20").Select()

$wb.Save()
